$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number need NumberFormat forced to
# Text ("@") first, otherwise Excel auto-converts the assigned string into a
# numeric value (dropping trailing zeros / changing type away from string).
$ws.Range('D2').Value = '26.054.85'
$ws.Range('E2').Value = '  +0.48%  '
$ws.Range('D3').Value = '1.639.57'
$ws.Range('E3').Value = '  +0.18%  '
$ws.Range('E4').Value = '  +0.58%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.71'
$ws.Range('E5').Value = '  -0.32%  '
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  -1.74%  '
$ws.Range('E9').Value = '  -1.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.68'
$ws.Range('E10').Value = '  -4.62%  '
$ws.Range('E11').Value = '  +0.04%  '
$ws.Range('B12').Value = 'Polkadot'
$ws.Range('C12').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.21'
$ws.Range('E12').Value = '  -1.48%  '
$ws.Range('B13').Value = 'WrappedEther'
$ws.Range('C13').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D13').Value = '1.622.00'
$ws.Range('E13').Value = '  +0.28%  '
$ws.Range('E14').Value = '  -2.04%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '62.41'
$ws.Range('E15').Value = '  -0.70%  '
$ws.Range('E16').Value = '  -1.87%  '
$ws.Range('D17').Value = '26.065.56'
$ws.Range('E17').Value = '  +0.59%  '
$ws.Range('E18').Value = '  +0.49%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '191.37'
$ws.Range('E19').Value = '  -0.63%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.27'
$ws.Range('E20').Value = '  -1.73%  '
$ws.Range('E21').Value = '  -2.94%  '
$ws.Range('E22').Value = '  -1.72%  '
$ws.Range('E23').Value = '  +0.32%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '144.07'
$ws.Range('E24').Value = '  +0.51%  '
$ws.Range('E25').Value = '  +0.64%  '
$ws.Range('E26').Value = '  -0.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '6.79'
$ws.Range('E27').Value = '  -1.32%  '
$ws.Range('E28').Value = '  -1.85%  '
$ws.Range('E29').Value = '  -0.20%  '
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.17'
$ws.Range('E31').Value = '  -2.18%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.18'
$ws.Range('E32').Value = '  -2.99%  '
$ws.Range('E33').Value = '  -1.02%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.43'
$ws.Range('E34').Value = '  -0.98%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.878'
$ws.Range('E35').Value = '  -2.45%  '
$ws.Range('D36').Value = '1.132.42'
$ws.Range('E36').Value = '  +0.10%  '
$ws.Range('E38').Value = '  -3.02%  '
$ws.Range('E39').Value = '  -1.06%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '98.92'
$ws.Range('E40').Value = '  -0.29%  '
$ws.Range('E41').Value = '  -1.22%  '
$ws.Range('E42').Value = '  -2.86%  '
$ws.Range('E43').Value = '  -0.90%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '55.42'
$ws.Range('E44').Value = '  -2.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0527'
$ws.Range('E45').Value = '  -0.31%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.49'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('E47').Value = '  +0.08%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.60'
$ws.Range('E48').Value = '  -0.77%  '
$ws.Range('E49').Value = '  +0.16%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0930'
$ws.Range('E50').Value = '  -2.93%  '
$ws.Range('E51').Value = '  -0.35%  '
